$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [double]"-6.988130962781725E-08"
$ws.Range("C3").Value = [double]"0"
$ws.Range("D3").Value = [double]"0"
$ws.Range("E3").Value = [double]"-0.05830520788542398"
$ws.Range("F3").Value = [double]"0.3400238586190717"
$ws.Range("G3").Value = [double]"0.07535685496909933"
$ws.Range("H3").Value = [double]"0.9063929398773607"
$ws.Range("I3").Value = [double]"0.7116691107469898"
$ws.Range("J3").Value = [double]"0.6200230953650364"
$ws.Range("K3").Value = [double]"2.427233269204584"
$ws.Range("L3").Value = [double]"0.00807063324930978"
$ws.Range("M3").Value = [double]"0.005430158360001892"
$ws.Range("C4").Value = [double]"0"
$ws.Range("D4").Value = [double]"0"
$ws.Range("E4").Value = [double]"-0.0147603227338363"
$ws.Range("F4").Value = [double]"0.4014853019561169"
$ws.Range("G4").Value = [double]"0.05908379740291336"
$ws.Range("H4").Value = [double]"1.008302110749463"
$ws.Range("I4").Value = [double]"0.5876954894759492"
$ws.Range("J4").Value = [double]"0.4783395297021388"
$ws.Range("K4").Value = [double]"2.355380386974529"
$ws.Range("L4").Value = [double]"0.007421882828360582"
$ws.Range("M4").Value = [double]"0.004900623612891733"
$ws.Range("C5").Value = [double]"0"
$ws.Range("D5").Value = [double]"0"
$ws.Range("E5").Value = [double]"-0.0147603227338363"
$ws.Range("F5").Value = [double]"0.4014853019561169"
$ws.Range("G5").Value = [double]"0.05908379740291336"
$ws.Range("H5").Value = [double]"1.008302110749463"
$ws.Range("I5").Value = [double]"0.5876954894759492"
$ws.Range("J5").Value = [double]"0.4783395297021388"
$ws.Range("K5").Value = [double]"2.355380386974529"
$ws.Range("L5").Value = [double]"0.007421882828360582"
$ws.Range("M5").Value = [double]"0.004900623612891733"
$ws.Range("B6").Value = [double]"0"
$ws.Range("C6").Value = [double]"0"
$ws.Range("D6").Value = [double]"0"
$ws.Range("E6").Value = [double]"0.04635286573504937"
$ws.Range("F6").Value = [double]"0.3934541319018796"
$ws.Range("G6").Value = [double]"-0.04984878472606311"
$ws.Range("H6").Value = [double]"0.4154546887941008"
$ws.Range("I6").Value = [double]"0.0007553802020709633"
$ws.Range("J6").Value = [double]"0.0002235581476850461"
$ws.Range("K6").Value = [double]"1.528949768037249"
$ws.Range("L6").Value = [double]"2.839899258795642E-29"
$ws.Range("M6").Value = [double]"2.845993463213659E-29"
$ws.Range("C7").Value = [double]"0"
$ws.Range("D7").Value = [double]"0"
$ws.Range("E7").Value = [double]"0.02968532041773266"
$ws.Range("F7").Value = [double]"0.3421199798815353"
$ws.Range("G7").Value = [double]"-0.05687896970063191"
$ws.Range("H7").Value = [double]"0.3504061492434056"
$ws.Range("I7").Value = [double]"0.0005230405822047129"
$ws.Range("J7").Value = [double]"0.000123524039266873"
$ws.Range("K7").Value = [double]"1.354212352525112"
$ws.Range("L7").Value = [double]"3.594247499413235E-29"
$ws.Range("M7").Value = [double]"3.601960476879787E-29"
$ws.Range("C8").Value = [double]"0"
$ws.Range("D8").Value = [double]"0"
$ws.Range("E8").Value = [double]"0.1480759987840615"
$ws.Range("F8").Value = [double]"0.7659614723941388"
$ws.Range("G8").Value = [double]"0.2935634188803161"
$ws.Range("H8").Value = [double]"0.3128696126156332"
$ws.Range("I8").Value = [double]"0.0004536009966637378"
$ws.Range("J8").Value = [double]"0.0002023950798711012"
$ws.Range("K8").Value = [double]"2.118143040633874"
$ws.Range("L8").Value = [double]"2.415886522239349E-28"
$ws.Range("M8").Value = [double]"2.421070828081064E-28"
$ws.Range("C9").Value = [double]"0"
$ws.Range("D9").Value = [double]"0"
$ws.Range("E9").Value = [double]"-0.003579967543327775"
$ws.Range("F9").Value = [double]"3.009209104201039"
$ws.Range("G9").Value = [double]"0.09493837353832454"
$ws.Range("H9").Value = [double]"0.2903324325689372"
$ws.Range("I9").Value = [double]"0.1791643622866484"
$ws.Range("J9").Value = [double]"0.01007169754559157"
$ws.Range("K9").Value = [double]"0.801944788792938"
$ws.Range("L9").Value = [double]"5.435744675038534E-30"
$ws.Range("M9").Value = [double]"4.940960873634824E-30"
$ws.Range("C10").Value = [double]"0"
$ws.Range("D10").Value = [double]"0"
$ws.Range("E10").Value = [double]"0.04752840069908072"
$ws.Range("F10").Value = [double]"4.817190569244339"
$ws.Range("G10").Value = [double]"-0.2654896072804981"
$ws.Range("H10").Value = [double]"0.484355037267985"
$ws.Range("I10").Value = [double]"0.41937787721498"
$ws.Range("J10").Value = [double]"-0.03034066380269484"
$ws.Range("K10").Value = [double]"0.801944788792938"
$ws.Range("L10").Value = [double]"5.435744675038534E-30"
$ws.Range("M10").Value = [double]"4.940960873634824E-30"
$ws.Range("C11").Value = [double]"0"
$ws.Range("D11").Value = [double]"0"
$ws.Range("E11").Value = [double]"0.04752840069908072"
$ws.Range("F11").Value = [double]"4.817190569244339"
$ws.Range("G11").Value = [double]"-0.2654896072804981"
$ws.Range("H11").Value = [double]"0.484355037267985"
$ws.Range("I11").Value = [double]"0.41937787721498"
$ws.Range("J11").Value = [double]"-0.03034066380269484"
$ws.Range("K11").Value = [double]"0.801944788792938"
$ws.Range("L11").Value = [double]"5.435744675038534E-30"
$ws.Range("M11").Value = [double]"4.940960873634824E-30"
$ws.Range("B12").Value = [double]"0"
$ws.Range("C12").Value = [double]"0"
$ws.Range("D12").Value = [double]"0"
$ws.Range("E12").Value = [double]"0.05978570338108916"
$ws.Range("F12").Value = [double]"2.101859155815864"
$ws.Range("G12").Value = [double]"-0.5222631437339688"
$ws.Range("H12").Value = [double]"1.923559800402738"
$ws.Range("I12").Value = [double]"0.01362372702717119"
$ws.Range("J12").Value = [double]"0.003307946537130625"
$ws.Range("K12").Value = [double]"1.282445377314577"
$ws.Range("L12").Value = [double]"5.965760595733902E-30"
$ws.Range("M12").Value = [double]"5.978562657098138E-30"
$ws.Range("C13").Value = [double]"0"
$ws.Range("D13").Value = [double]"0"
$ws.Range("E13").Value = [double]"0.1036930408465265"
$ws.Range("F13").Value = [double]"1.641608615122233"
$ws.Range("G13").Value = [double]"-0.0216343602596395"
$ws.Range("H13").Value = [double]"1.88335766327739"
$ws.Range("I13").Value = [double]"0.01377767394134808"
$ws.Range("J13").Value = [double]"0.0006950414970530134"
$ws.Range("K13").Value = [double]"2.034950389185837"
$ws.Range("L13").Value = [double]"1.774937036747277E-30"
$ws.Range("M13").Value = [double]"1.778745914508537E-30"
$ws.Range("C14").Value = [double]"0"
$ws.Range("D14").Value = [double]"0"
$ws.Range("E14").Value = [double]"0.1513245328683536"
$ws.Range("F14").Value = [double]"1.430229496766241"
$ws.Range("G14").Value = [double]"0.2038468732994141"
$ws.Range("H14").Value = [double]"1.20288124391921"
$ws.Range("I14").Value = [double]"0.005939583161199985"
$ws.Range("J14").Value = [double]"0.0009342452166989785"
$ws.Range("K14").Value = [double]"2.342500249112021"
$ws.Range("L14").Value = [double]"7.119469669619632E-29"
$ws.Range("M14").Value = [double]"7.905537397815719E-29"
$ws.Range("E15").Value = [double]"3.04982991176699E-05"
$ws.Range("F15").Value = [double]"1.322257755762889E-06"
$ws.Range("G15").Value = [double]"1.687078097113069E-07"
$ws.Range("H15").Value = [double]"2.853234980418634E-06"
$ws.Range("I15").Value = [double]"1.155619596986963E-11"
$ws.Range("J15").Value = [double]"1.049031921478155E-11"
$ws.Range("K15").Value = [double]"5.735841689583602E-06"
$ws.Range("L15").Value = [double]"2.524842922867389E-11"
$ws.Range("M15").Value = [double]"2.519518469257158E-11"
$ws.Range("E16").Value = [double]"-0.0002261597385029233"
$ws.Range("F16").Value = [double]"8.853544808984725E-06"
$ws.Range("G16").Value = [double]"1.816703049353086E-06"
$ws.Range("H16").Value = [double]"2.282970355035119E-05"
$ws.Range("I16").Value = [double]"2.126876298157896E-09"
$ws.Range("J16").Value = [double]"1.964295578190866E-09"
$ws.Range("K16").Value = [double]"4.138219880126974E-05"
$ws.Range("L16").Value = [double]"2.044025140945604E-09"
$ws.Range("M16").Value = [double]"2.041774794697484E-09"
$ws.Range("E17").Value = [double]"-0.0002465707462691049"
$ws.Range("F17").Value = [double]"9.347951400936354E-06"
$ws.Range("G17").Value = [double]"2.30918997797262E-06"
$ws.Range("H17").Value = [double]"2.371805848418646E-05"
$ws.Range("I17").Value = [double]"2.333542980756478E-09"
$ws.Range("J17").Value = [double]"2.166827275844826E-09"
$ws.Range("K17").Value = [double]"4.236582307273997E-05"
$ws.Range("L17").Value = [double]"2.134699070570012E-09"
$ws.Range("M17").Value = [double]"2.132903609809513E-09"
$ws.Range("E18").Value = [double]"0.02017975204408363"
$ws.Range("F18").Value = [double]"0.0004542608191382014"
$ws.Range("G18").Value = [double]"0.0004408710129724271"
$ws.Range("H18").Value = [double]"1.237277485020459E-05"
$ws.Range("I18").Value = [double]"1.500328420198403E-10"
$ws.Range("J18").Value = [double]"1.502598105907721E-10"
$ws.Range("K18").Value = [double]"0.001342053523274252"
$ws.Range("L18").Value = [double]"9.4544502669448E-07"
$ws.Range("M18").Value = [double]"9.466131698347057E-07"
$ws.Range("E19").Value = [double]"0.006143124466848736"
$ws.Range("F19").Value = [double]"5.699515167131308E-05"
$ws.Range("G19").Value = [double]"5.595884323394623E-05"
$ws.Range("H19").Value = [double]"1.078460907728614E-06"
$ws.Range("I19").Value = [double]"7.203138133780777E-13"
$ws.Range("J19").Value = [double]"7.199294056447268E-13"
$ws.Range("K19").Value = [double]"0.0002574414881415026"
$ws.Range("L19").Value = [double]"6.873472437492059E-09"
$ws.Range("M19").Value = [double]"6.848251690458709E-09"
$ws.Range("E20").Value = [double]"0.006143124466848736"
$ws.Range("F20").Value = [double]"5.699515167131308E-05"
$ws.Range("G20").Value = [double]"5.595884323394623E-05"
$ws.Range("H20").Value = [double]"1.078460907728614E-06"
$ws.Range("I20").Value = [double]"7.203138133780777E-13"
$ws.Range("J20").Value = [double]"7.199294056447268E-13"
$ws.Range("K20").Value = [double]"0.0002574414881415026"
$ws.Range("L20").Value = [double]"6.873472437492059E-09"
$ws.Range("M20").Value = [double]"6.848251690458709E-09"
$ws.Range("E21").Value = [double]"-9.550550007629545E-06"
$ws.Range("F21").Value = [double]"2.322636844196244E-06"
$ws.Range("G21").Value = [double]"8.49110860072011E-07"
$ws.Range("H21").Value = [double]"0.001234793694234009"
$ws.Range("I21").Value = [double]"4.293935691277007E-06"
$ws.Range("J21").Value = [double]"4.355340793075507E-07"
$ws.Range("K21").Value = [double]"9.861348368513163E-07"
$ws.Range("L21").Value = [double]"5.271151386122459E-13"
$ws.Range("M21").Value = [double]"5.038616579524655E-13"
$ws.Range("E22").Value = [double]"-4.918013016395333E-05"
$ws.Range("F22").Value = [double]"1.212223120551698E-06"
$ws.Range("G22").Value = [double]"-3.434865358340723E-08"
$ws.Range("H22").Value = [double]"3.918925665006828E-06"
$ws.Range("I22").Value = [double]"5.198377189635467E-11"
$ws.Range("J22").Value = [double]"1.089940528447007E-11"
$ws.Range("K22").Value = [double]"2.015933335195811E-06"
$ws.Range("L22").Value = [double]"4.060592043585097E-12"
$ws.Range("M22").Value = [double]"3.945789335512037E-12"
$ws.Range("E23").Value = [double]"-4.918013016395333E-05"
$ws.Range("F23").Value = [double]"1.212223120551698E-06"
$ws.Range("G23").Value = [double]"-3.434865358340723E-08"
$ws.Range("H23").Value = [double]"3.918925665006828E-06"
$ws.Range("I23").Value = [double]"5.198377189635467E-11"
$ws.Range("J23").Value = [double]"1.089940528447007E-11"
$ws.Range("K23").Value = [double]"2.015933335195811E-06"
$ws.Range("L23").Value = [double]"4.060592043585097E-12"
$ws.Range("M23").Value = [double]"3.945789335512037E-12"
$ws.Range("E24").Value = [double]"0.02313835960025783"
$ws.Range("F24").Value = [double]"0.0004687964129736069"
$ws.Range("G24").Value = [double]"0.000406646230379195"
$ws.Range("H24").Value = [double]"9.871918216865349E-06"
$ws.Range("I24").Value = [double]"3.908704261820706E-11"
$ws.Range("J24").Value = [double]"3.81351275277919E-11"
$ws.Range("K24").Value = [double]"0.0005655587627048674"
$ws.Range("L24").Value = [double]"1.118670363394964E-07"
$ws.Range("M24").Value = [double]"1.091460264437043E-07"
$ws.Range("E25").Value = [double]"0.02458316731348448"
$ws.Range("F25").Value = [double]"0.000357009619959635"
$ws.Range("G25").Value = [double]"0.0003104980422480224"
$ws.Range("H25").Value = [double]"3.63863074993663E-05"
$ws.Range("I25").Value = [double]"7.021105256404883E-10"
$ws.Range("J25").Value = [double]"6.863756048990907E-10"
$ws.Range("K25").Value = [double]"0.0004601043896032899"
$ws.Range("L25").Value = [double]"1.060248579879899E-07"
$ws.Range("M25").Value = [double]"1.036819053587068E-07"
$ws.Range("E26").Value = [double]"0.1411814472911681"
$ws.Range("F26").Value = [double]"0.007451042704150982"
$ws.Range("G26").Value = [double]"0.006512965616889736"
$ws.Range("H26").Value = [double]"0.02460215372691189"
$ws.Range("I26").Value = [double]"0.0004235135650665852"
$ws.Range("J26").Value = [double]"0.0004152128383758606"
$ws.Range("K26").Value = [double]"0.02600202288483916"
$ws.Range("L26").Value = [double]"0.000520414257377082"
$ws.Range("M26").Value = [double]"0.0006010811951497761"
